$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the AO3:AO18 literal values (per commit message: "sign problem" fix,
# equilibrium equation constant updated from 611250.58749924635 to 578039.49484035082)
$newValue = 578039.49484035082

for ($row = 3; $row -le 18; $row++) {
    $ws.Range("AO$row").Value = $newValue
}

$excel.Calculate()
